$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-38: refreshed price (D) and 1h volume change (E) values
# (leading apostrophe forces numeric-looking price strings to stay text,
#  matching the original inlineStr cell type)
$ws.Range("D2").Value = "66.327.14"
$ws.Range("E2").Value = "  +6.57%  "
$ws.Range("D3").Value = "3.005.04"
$ws.Range("E3").Value = "  +3.41%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'582.16"
$ws.Range("E5").Value = "  +2.65%  "
$ws.Range("D6").Value = "'162.14"
$ws.Range("E6").Value = "  +12.59%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "  +3.82%  "
$ws.Range("D9").Value = "3.001.57"
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("D10").Value = "'6.59"
$ws.Range("E10").Value = "  -5.04%  "
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("D12").Value = "'0.456"
$ws.Range("E12").Value = "  +5.24%  "
$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  +6.59%  "
$ws.Range("D14").Value = "'34.60"
$ws.Range("E14").Value = "  +6.04%  "
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "66.250.07"
$ws.Range("E16").Value = "  +6.62%  "
$ws.Range("D17").Value = "3.502.85"
$ws.Range("E17").Value = "  +3.45%  "
$ws.Range("D18").Value = "'6.91"
$ws.Range("E18").Value = "  +5.34%  "
$ws.Range("D19").Value = "3.004.08"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("D20").Value = "'455.79"
$ws.Range("E20").Value = "  +6.37%  "
$ws.Range("D21").Value = "'13.82"
$ws.Range("E21").Value = "  +5.81%  "
$ws.Range("E22").Value = "  +4.45%  "
$ws.Range("D23").Value = "'7.34"
$ws.Range("E23").Value = "  +7.32%  "
$ws.Range("D24").Value = "'82.28"
$ws.Range("E24").Value = "  +4.77%  "
$ws.Range("D25").Value = "'2.30"
$ws.Range("E25").Value = "  +15.08%  "
$ws.Range("D26").Value = "'12.30"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("D27").Value = "'10.46"
$ws.Range("E27").Value = "  +5.45%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'8.19"
$ws.Range("E29").Value = "  +18.16%  "
$ws.Range("E30").Value = "  +20.72%  "
$ws.Range("E31").Value = "  -4.67%  "
$ws.Range("D32").Value = "'2.61"
$ws.Range("E32").Value = "  +5.19%  "
$ws.Range("D33").Value = "'27.20"
$ws.Range("E33").Value = "  +6.23%  "
$ws.Range("D34").Value = "'0.111"
$ws.Range("E34").Value = "  +4.68%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'0.993"
$ws.Range("E36").Value = "  +4.45%  "
$ws.Range("D37").Value = "'5.80"
$ws.Range("E37").Value = "  +7.88%  "
$ws.Range("D38").Value = "'2.16"
$ws.Range("E38").Value = "  +14.39%  "

# Rows 39-40: OKB and dogwifhat swapped list positions, with refreshed price/volume
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  +2.47%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'49.89"
$ws.Range("E40").Value = "  +2.11%  "

# Rows 41-51: refreshed price (D) and 1h volume change (E) values
$ws.Range("D41").Value = "'0.310"
$ws.Range("E41").Value = "  +16.95%  "
$ws.Range("E42").Value = "  +7.55%  "
$ws.Range("D43").Value = "'43.99"
$ws.Range("E43").Value = "  +7.40%  "
$ws.Range("D44").Value = "'8.41"
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("D45").Value = "'398.53"
$ws.Range("E45").Value = "  +14.71%  "
$ws.Range("E46").Value = "  +7.48%  "
$ws.Range("D47").Value = "2.794.28"
$ws.Range("E47").Value = "  +3.36%  "
$ws.Range("D48").Value = "'134.76"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'23.91"
$ws.Range("E50").Value = "  +11.93%  "
$ws.Range("E51").Value = "  +4.80%  "
